$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new history rows (27 and 28) to the project history log, mirroring
# the formatting of the existing last row (26).
# ---------------------------------------------------------------------------

# Row 27: copy formatting (without column E) from row 26, then fill values.
$ws.Range("A26:D26").Copy()
$ws.Range("A27:D27").PasteSpecial(-4122)

# Row 28: copy formatting (including column E) from row 26, then fill values.
$ws.Range("A26:E26").Copy()
$ws.Range("A28:E28").PasteSpecial(-4122)

# --- Row 27 values ----------------------------------------------------------
$ws.Range("A27").Value = 44987
$ws.Range("B27").Value = 24
$ws.Range("C27").Value = "Switch to the new branch by using 'git switch -c ""revised/side-menu-lists""'"
$ws.Range("D27").Value = "git switch -c ""revised/side-menu-lists"""

# --- Row 28 values ----------------------------------------------------------
$ws.Range("A28").Value = 44987
$ws.Range("B28").Value = 25
$ws.Range("C28").Value = "Revise side menu and routing"
$ws.Range("D28").Value = "Menu`n1. Home`n2. My performance`n3. SEP Card`n4. E-TS1`n   4.1 My Work space`n   4.2 My E-TS1"
$ws.Range("E28").Value = "pages-menu.ts`npages-routing.module.ts"

# Row 28 wraps across 7 lines in column D, so Excel grows the row to 105pt.
$ws.Rows.Item(28).RowHeight = 105

# ---------------------------------------------------------------------------
# Grow the Table2 structured table so it covers the two new rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E28"))

# ---------------------------------------------------------------------------
# Update the view so the new rows are visible, matching the saved selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("A29").Select()
